# Weekly Fruta/Hortaliza update: insert 3 new daily records (2021-10-05,
# "Cultivar IV Región" / "Provincia del Elquí") ahead of the existing
# Chirimoya rows, pushing the rest of the table down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 142, shifting rows 142:159 down to 145:162.
$ws.Rows("142:144").Insert()

# Values shared by every row in this sub-table.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "`$/kilo (en caja de 15 kilos)"
$origen      = "Provincia del Elquí"
$fecha       = 44491

$newRows = @(
    @{ Row = 142; Calidad = "Especial";                Volumen = 150; Precio = 2700 },
    @{ Row = 143; Calidad = "Extra (doble especial)";  Volumen = 75;  Precio = 2900 },
    @{ Row = 144; Calidad = "Primera";                 Volumen = 125; Precio = 2500 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $mercadoId
    $ws.Range("B$row").Value2 = $mercado
    $ws.Range("C$row").Value2 = $region
    $ws.Range("D$row").Value2 = $fecha
    $ws.Range("E$row").Value2 = $codreg
    $ws.Range("F$row").Value2 = $tipo
    $ws.Range("G$row").Value2 = $productoId
    $ws.Range("H$row").Value2 = $producto
    $ws.Range("I$row").Value2 = $categoriaId
    $ws.Range("J$row").Value2 = $categoria
    $ws.Range("K$row").Value2 = $variedad
    $ws.Range("L$row").Value2 = $r.Calidad
    $ws.Range("M$row").Value2 = $r.Volumen
    $ws.Range("N$row").Value2 = $r.Precio
    $ws.Range("O$row").Value2 = $r.Precio
    $ws.Range("P$row").Value2 = $r.Precio
    $ws.Range("Q$row").Value2 = $unidad
    $ws.Range("R$row").Value2 = $origen
    $ws.Range("S$row").Value2 = $r.Precio
    $ws.Range("T$row").Value2 = 1
}
